# country_parameters.xlsx update:
#  - replace the two sample countries (Kenya / Tanzania) with Namibia / Other
#  - refresh the interest-rate / price figures to the new shared value
#    (5.7890216323739932E-2) used across solar / wind / plant / infrastructure,
#    plus new electricity & heat prices
#  - give the "Infrastructure interest rate" column (J) a high-precision
#    number format so the new shared rate isn't rounded for display

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Country names (column A) --------------------------------------------
$ws.Range("A2").Value = "Namibia"
$ws.Range("A3").Value = "Other"

# -- Row 2 (Namibia) -------------------------------------------------------
$ws.Range("B2").Value = 0.10465000000000001   # Electricity price (euros/kWh)
$ws.Range("C2").Value = 0.02                  # Heat price (euros/kWh)
$ws.Range("D2").Value = 0.057890216323739932  # Solar interest rate
$ws.Range("E2").Value = 20                    # Solar lifetime (years)
$ws.Range("F2").Value = 0.057890216323739932  # Wind interest rate
$ws.Range("G2").Value = 20                    # Wind lifetime (years)
$ws.Range("H2").Value = 0.057890216323739932  # Plant interest rate
$ws.Range("I2").Value = 20                    # Plant lifetime (years)
$ws.Range("J2").Value = 0.057890216323739932  # Infrastructure interest rate
$ws.Range("K2").Value = 50                    # Infrastructure lifetime (years)

# -- Row 3 (Other) ----------------------------------------------------------
$ws.Range("B3").Value = 0.10465000000000001
$ws.Range("C3").Value = 0.02
$ws.Range("D3").Value = 0.057890216323739932
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 0.057890216323739932
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 0.057890216323739932
$ws.Range("I3").Value = 20
$ws.Range("J3").Value = 0.057890216323739932
$ws.Range("K3").Value = 50

# -- Number format for the infrastructure-interest-rate column -------------
$ws.Range("J2:J3").NumberFormat = "0.000000000"

# -- Cursor/selection where the author last left it -------------------------
$ws.Range("A4").Select()
